$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: add PriceChange (X8) and UpDown (Y8)
$ws.Range("X8").Value = 0.29000100000000373
$ws.Range("Y8").Value = "Up"

# Row 9: new data row
$ws.Range("A9").Value = 42648.890601851854
$ws.Range("B9").Value = -9
$ws.Range("C9").Value = "Sell"
$ws.Range("D9").Value = 20
$ws.Range("E9").Value = 24682
$ws.Range("F9").Value = 2950
$ws.Range("G9").Value = 58
$ws.Range("H9").Value = 37
$ws.Range("I9").Value = 81
$ws.Range("J9").Value = 17
$ws.Range("K9").Value = 34045
$ws.Range("L9").Value = 344
$ws.Range("M9").Value = 221
$ws.Range("N9").Value = 125
$ws.Range("O9").Value = 27
$ws.Range("P9").Value = "Bag"
$ws.Range("Q9").Value = 38.916275631518758
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = -0.0077999999999999996
$ws.Range("S9").NumberFormat = "0.00%"
$ws.Range("T9").Value = -0.030499999999999999
$ws.Range("T9").NumberFormat = "0.00%"
$ws.Range("U9").Value = 14.62
$ws.Range("V9").Value = "N/A"
$ws.Range("W9").Value = -2
